$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.870.76'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '1.893.82'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7830'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.20'
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3150'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.48'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07336'
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08130'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7683'
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.481'
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").Value = '1.879.42'
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.33'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.215'
$ws.Range("D17").Value = '29.857.26'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.97'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.89'
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.160'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '2.122.45'
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9986'
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.474'
$ws.Range("E26").Value = '  +1.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.17'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.81'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.040'
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.461'
$ws.Range("E30").Value = '  +6.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.544'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.483'
$ws.Range("E32").Value = '  +1.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05606'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.087'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.255'
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7570'
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.644'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01939'
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.781'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = '1.145.15'
$ws.Range("E41").Value = '  +11.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4466'
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.00'
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.970'
$ws.Range("E44").Value = '  +2.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8580'
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.904'
$ws.Range("E46").Value = '  +2.77%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9995'
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.124'
$ws.Range("E48").Value = '  +6.98%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.97'
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.827'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.531'
$ws.Range("E51").Value = '  +1.19%  '
